$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fixed Weight")

# The "color" column (D) listed a couple of fasteners as simply "zinc"
# plated hardware; update the description to the fuller "Zinc Plated"
# wording used elsewhere in the BOM.
$ws.Range("D6").Value = "Zinc Plated"
$ws.Range("D14").Value = "Zinc Plated"

# Autofit the color column now that its text is longer, and land the
# selection back on the header cell like the saved workbook shows.
$ws.Columns.Item(4).AutoFit() | Out-Null
$ws.Range("D2").Select() | Out-Null
